$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 0.0529
$ws.Range("F3").Value = 1292026
$ws.Range("C3").Select()

$wb.Save()
